$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3757192674132455
$ws.Range("C2").Value = 0.02930152855844881
$ws.Range("D2").Value = 0.2420222888857637
$ws.Range("F2").Value = 1.586398145579913
$ws.Range("G2").Value = 0.002472051365649399
$ws.Range("J2").Value = 0.2973570291902163
$ws.Range("K2").Value = 0.328076540996932
$ws.Range("M2").Value = 0.2898795205865525
$ws.Range("N2").Value = 1.891113143744201
$ws.Range("O2").Value = 3.731092495176853
$ws.Range("B3").Value = 0.3437020739106345
$ws.Range("C3").Value = 0.02575381342566629
$ws.Range("D3").Value = 0.2372607916900762
$ws.Range("F3").Value = 1.588810574898332
$ws.Range("G3").Value = 0.002474450459123423
$ws.Range("J3").Value = 0.2955377707892382
$ws.Range("K3").Value = 0.2951558770172937
$ws.Range("M3").Value = 0.2772209927184832
$ws.Range("N3").Value = 1.909769330211154
$ws.Range("O3").Value = 3.746912369428998
$ws.Range("B4").Value = 0.3241446689584393
$ws.Range("C4").Value = 0.02356485319731405
$ws.Range("D4").Value = 0.2344398604698057
$ws.Range("F4").Value = 1.591045768473869
$ws.Range("G4").Value = 0.00247600309385297
$ws.Range("J4").Value = 0.2945670943208114
$ws.Range("K4").Value = 0.2749924140743332
$ws.Range("M4").Value = 0.2695754061435309
$ws.Range("N4").Value = 1.921808005209301
$ws.Range("O4").Value = 3.758406694073443
$ws.Range("B5").Value = 0.3162007860394738
$ws.Range("C5").Value = 0.02267020060249081
$ws.Range("D5").Value = 0.2333162290993585
$ws.Range("F5").Value = 1.59214637645308
$ws.Range("G5").Value = 0.002476655875273182
$ws.Range("J5").Value = 0.2942083735783783
$ws.Range("K5").Value = 0.2667886434580993
$ws.Range("M5").Value = 0.2664918361548061
$ws.Range("N5").Value = 1.926860653544052
$ws.Range("O5").Value = 3.763538740710402
$ws.Range("B6").Value = 0.3148832890013011
$ws.Range("C6").Value = 0.02252148646562091
$ws.Range("D6").Value = 0.2331312197529201
$ws.Range("F6").Value = 1.592340596659724
$ws.Range("G6").Value = 0.002476765482952777
$ws.Range("J6").Value = 0.2941510342330602
$ws.Range("K6").Value = 0.265427212037153
$ws.Range("M6").Value = 0.2659817539759786
$ws.Range("N6").Value = 1.927708508193244
$ws.Range("O6").Value = 3.764417978278843
$ws.Range("B7").Value = 0.3240374294363164
$ws.Range("C7").Value = 0.02355279820453404
$ws.Range("D7").Value = 0.2344246016852907
$ws.Range("F7").Value = 1.591059843161013
$ws.Range("G7").Value = 0.002476011815953914
$ws.Range("J7").Value = 0.2945621072848894
$ws.Range("K7").Value = 0.2748817219152357
$ws.Range("M7").Value = 0.2695336899521408
$ws.Range("N7").Value = 1.921875552624334
$ws.Range("O7").Value = 3.758474092380766
$ws.Range("B8").Value = 0.3646589781525904
$ws.Range("C8").Value = 0.02808051201088801
$ws.Range("D8").Value = 0.2403592839521167
$ws.Range("F8").Value = 1.587073561752497
$ws.Range("G8").Value = 0.002472862091322307
$ws.Range("J8").Value = 0.2966993928219139
$ws.Range("K8").Value = 0.3167153821872262
$ws.Range("M8").Value = 0.2854886483864902
$ws.Range("N8").Value = 1.897424659926056
$ws.Range("O8").Value = 3.736177688758573
$ws.Range("B9").Value = 0.4451056395042201
$ws.Range("C9").Value = 0.03687331581217279
$ws.Range("D9").Value = 0.2528076325190369
$ws.Range("F9").Value = 1.585233192396188
$ws.Range("G9").Value = 0.002467314249548926
$ws.Range("J9").Value = 0.3020508801079558
$ws.Range("K9").Value = 0.3991321295707166
$ws.Range("M9").Value = 0.3177761352391713
$ws.Range("N9").Value = 1.854106567357467
$ws.Range("O9").Value = 3.706578236069276
$ws.Range("B10").Value = 0.5046750870809262
$ws.Range("C10").Value = 0.04327947746035932
$ws.Range("D10").Value = 0.2624433132386059
$ws.Range("F10").Value = 1.587519419146417
$ws.Range("G10").Value = 0.002463617796766438
$ws.Range("J10").Value = 0.3066894704071075
$ws.Range("K10").Value = 0.4599014553843404
$ws.Range("M10").Value = 0.3421015858874554
$ws.Range("N10").Value = 1.825100404069672
$ws.Range("O10").Value = 3.693436459953773
$ws.Range("B11").Value = 0.5318729635109207
$ws.Range("C11").Value = 0.04618182493973677
$ws.Range("D11").Value = 0.2669323866935827
$ws.Range("F11").Value = 1.589348569943468
$ws.Range("G11").Value = 0.002462017800038155
$ws.Range("J11").Value = 0.3089531229355487
$ws.Range("K11").Value = 0.4875915403929127
$ws.Range("M11").Value = 0.3532979264869169
$ws.Range("N11").Value = 1.81251637930111
$ws.Range("O11").Value = 3.689325586528696
$ws.Range("B12").Value = 0.5421860114604158
$ws.Range("C12").Value = 0.04727912763763698
$ws.Range("D12").Value = 0.2686473941411975
$ws.Range("F12").Value = 1.590154580148592
$ws.Range("G12").Value = 0.002461423587654051
$ws.Range("J12").Value = 0.309832365434886
$ws.Range("K12").Value = 0.4980832725769915
$ws.Range("M12").Value = 0.3575563150234018
$ws.Range("N12").Value = 1.807838993990059
$ws.Range("O12").Value = 3.688037315218565
$ws.Range("B13").Value = 0.5399643061044515
$ws.Range("C13").Value = 0.0470428826366458
$ws.Range("D13").Value = 0.2682773671737522
$ws.Range("F13").Value = 1.589975952146048
$ws.Range("G13").Value = 0.002461551043562765
$ws.Range("J13").Value = 0.3096420248173501
$ws.Range("K13").Value = 0.4958234269299453
$ws.Range("M13").Value = 0.3566383717956398
$ws.Range("N13").Value = 1.808842441926616
$ws.Range("O13").Value = 3.688302830720943
$ws.Range("B14").Value = 0.532721150396668
$ws.Range("C14").Value = 0.04627213611684056
$ws.Range("D14").Value = 0.2670731795643064
$ws.Range("F14").Value = 1.589412609898545
$ws.Range("G14").Value = 0.002461968680199719
$ws.Range("J14").Value = 0.3090250170842239
$ws.Range("K14").Value = 0.4884545821107338
$ws.Range("M14").Value = 0.3536478954966
$ws.Range("N14").Value = 1.812129805834168
$ws.Range("O14").Value = 3.689214220566527
$ws.Range("B15").Value = 0.5282862911243171
$ws.Range("C15").Value = 0.0457998018143968
$ws.Range("D15").Value = 0.2663375420016934
$ws.Range("F15").Value = 1.589082304339783
$ws.Range("G15").Value = 0.002462226013056542
$ws.Range("J15").Value = 0.3086499521824351
$ws.Range("K15").Value = 0.483941731468633
$ws.Range("M15").Value = 0.3518185558299578
$ws.Range("N15").Value = 1.814154861300856
$ws.Range("O15").Value = 3.689807427330749
$ws.Range("B16").Value = 0.5028996003825625
$ws.Range("C16").Value = 0.04308955957992566
$ws.Range("D16").Value = 0.2621520599861782
$ws.Range("F16").Value = 1.587415748471969
$ws.Range("G16").Value = 0.002463723996446204
$ws.Range("J16").Value = 0.3065446221963555
$ws.Range("K16").Value = 0.4580927232212559
$ws.Range("M16").Value = 0.341372489155475
$ws.Range("N16").Value = 1.825935099028856
$ws.Range("O16").Value = 3.693742683176595
$ws.Range("B17").Value = 0.4873508028592539
$ws.Range("C17").Value = 0.04142384397347598
$ws.Range("D17").Value = 0.2596114070408788
$ws.Range("F17").Value = 1.586595397092879
$ws.Range("G17").Value = 0.002464663806500723
$ws.Range("J17").Value = 0.3052923718583997
$ws.Range("K17").Value = 0.4422465794235109
$ws.Range("M17").Value = 0.3349974668311262
$ws.Range("N17").Value = 1.83331842404554
$ws.Range("O17").Value = 3.696635045914888
$ws.Range("B18").Value = 0.4784169348648959
$ws.Range("C18").Value = 0.04046465622596429
$ws.Range("D18").Value = 0.2581600515916165
$ws.Range("F18").Value = 1.58619783732388
$ws.Range("G18").Value = 0.002465212038413666
$ws.Range("J18").Value = 0.3045865638900267
$ws.Range("K18").Value = 0.4331366381674684
$ws.Range("M18").Value = 0.3313430244795512
$ws.Range("N18").Value = 1.837622626819927
$ws.Range("O18").Value = 3.698474429690634
$ws.Range("B19").Value = 0.4753937083664823
$ws.Range("C19").Value = 0.04013970223009267
$ws.Range("D19").Value = 0.2576703618585015
$ws.Range("F19").Value = 1.586075992719401
$ws.Range("G19").Value = 0.002465398980425063
$ws.Range("J19").Value = 0.3043500727700064
$ws.Range("K19").Value = 0.43005292917681
$ws.Range("M19").Value = 0.330107811993706
$ws.Range("N19").Value = 1.839089831845981
$ws.Range("O19").Value = 3.699127405733776
$ws.Range("B20").Value = 0.4890050307568856
$ws.Range("C20").Value = 0.04160127771046973
$ws.Range("D20").Value = 0.2598808337949521
$ws.Range("F20").Value = 1.586675037627373
$ws.Range("G20").Value = 0.002464562968021222
$ws.Range("J20").Value = 0.3054241804106255
$ws.Range("K20").Value = 0.4439329826970777
$ws.Range("M20").Value = 0.3356748270801191
$ws.Range("N20").Value = 1.832526504139247
$ws.Range("O20").Value = 3.696308957819781
$ws.Range("B21").Value = 0.5348482701696753
$ws.Range("C21").Value = 0.04649857103295574
$ws.Range("D21").Value = 0.2674264699557369
$ws.Range("F21").Value = 1.589575001878458
$ws.Range("G21").Value = 0.002461845693847483
$ws.Range("J21").Value = 0.3092056491320676
$ws.Range("K21").Value = 0.4906188288685485
$ws.Range("M21").Value = 0.3545257682776324
$ws.Range("N21").Value = 1.811161841078292
$ws.Range("O21").Value = 3.688939239038689
$ws.Range("B22").Value = 0.5648896974727222
$ws.Range("C22").Value = 0.04968899271192129
$ws.Range("D22").Value = 0.2724458982068398
$ws.Range("F22").Value = 1.592130930933649
$ws.Range("G22").Value = 0.002460137804858532
$ws.Range("J22").Value = 0.3118055460638089
$ws.Range("K22").Value = 0.5211660439026389
$ws.Range("M22").Value = 0.3669541385168955
$ws.Range("N22").Value = 1.797711239535181
$ws.Range("O22").Value = 3.685687254678015
$ws.Range("B23").Value = 0.5488488463540477
$ws.Range("C23").Value = 0.04798715820857069
$ws.Range("D23").Value = 0.2697589282075796
$ws.Range("F23").Value = 1.590706377090171
$ws.Range("G23").Value = 0.002461043132331864
$ws.Range("J23").Value = 0.3104061859513791
$ws.Range("K23").Value = 0.5048593509829402
$ws.Range("M23").Value = 0.3603110497484181
$ws.Range("N23").Value = 1.80484317406336
$ws.Range("O23").Value = 3.687279776417711
$ws.Range("B24").Value = 0.4882571376452347
$ws.Range("C24").Value = 0.04152106472947992
$ws.Range("D24").Value = 0.2597589969841607
$ws.Range("F24").Value = 1.586638801404774
$ws.Range("G24").Value = 0.002464608532313951
$ws.Range("J24").Value = 0.305364545756845
$ws.Range("K24").Value = 0.4431705590617696
$ws.Range("M24").Value = 0.3353685593458522
$ws.Range("N24").Value = 1.832884346005517
$ws.Range("O24").Value = 3.69645583236877
$ws.Range("B25").Value = 0.4232598877458429
$ws.Range("C25").Value = 0.03450399110559488
$ws.Range("D25").Value = 0.2493536927101303
$ws.Range("F25").Value = 1.58509175246796
$ws.Range("G25").Value = 0.002468748165417864
$ws.Range("J25").Value = 0.3004789947328632
$ws.Range("K25").Value = 0.3767969246239318
$ws.Range("M25").Value = 0.3089350609703772
$ws.Range("N25").Value = 1.865329898336183
$ws.Range("O25").Value = 3.713074177729084
